# Apply the "Updated cryptos list" refresh: new prices / volume percentages,
# plus a few rows whose coin (name/link) swapped rank position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new text value.
# Numeric-looking values are prefixed with a literal leading apostrophe so
# Excel stores them as text (matching the sheet existing string columns)
# instead of silently converting them to floating-point numbers.
$updates = [ordered]@{
    "D2" = "31.384.64"
    "E2" = "  +3.58%  "
    "D3" = "2.006.37"
    "E3" = "  +7.46%  "
    "D4" = "'0.9993"
    "E4" = "  -0.09%  "
    "D5" = "'0.7937"
    "E5" = "  +67.94%  "
    "D6" = "'259.77"
    "E6" = "  +6.68%  "
    "D7" = "'0.9983"
    "E7" = "  -0.17%  "
    "D8" = "'0.3608"
    "E8" = "  +25.57%  "
    "D9" = "'28.66"
    "E9" = "  +33.16%  "
    "D10" = "'0.07056"
    "E10" = "  +8.93%  "
    "D11" = "'0.8464"
    "E11" = "  +18.00%  "
    "D12" = "'0.08089"
    "E12" = "  +3.91%  "
    "B13" = "Litecoin"
    "C13" = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
    "D13" = "'101.35"
    "E13" = "  +4.88%  "
    "B14" = "WrappedEther"
    "C14" = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
    "D14" = "2.002.16"
    "E14" = "  +7.24%  "
    "D15" = "'5.636"
    "E15" = "  +9.98%  "
    "D16" = "'275.80"
    "E16" = "  -1.74%  "
    "D17" = "31.359.77"
    "E17" = "  +3.54%  "
    "D18" = "'14.66"
    "E18" = "  +12.80%  "
    "D19" = "'5.928"
    "E19" = "  +13.13%  "
    "D20" = "'0.000007932"
    "E20" = "  +6.26%  "
    "D21" = "2.263.71"
    "D22" = "'0.9975"
    "E22" = "  -0.24%  "
    "D23" = "'0.9995"
    "E23" = "  -0.08%  "
    "D24" = "'7.201"
    "E24" = "  +15.19%  "
    "D25" = "'10.22"
    "E25" = "  +13.89%  "
    "D26" = "'0.1522"
    "E26" = "  +57.98%  "
    "D27" = "'164.47"
    "E27" = "  +1.47%  "
    "D28" = "'20.06"
    "E28" = "  +7.41%  "
    "D29" = "'2.408"
    "E29" = "  +28.28%  "
    "D30" = "'1.634"
    "E30" = "  +10.24%  "
    "D31" = "'4.645"
    "E31" = "  +10.41%  "
    "D32" = "'1.357"
    "E32" = "  +3.55%  "
    "D33" = "'4.407"
    "E33" = "  +7.15%  "
    "D34" = "'0.05225"
    "E34" = "  +9.31%  "
    "D35" = "'1.224"
    "E35" = "  +9.65%  "
    "D36" = "'0.7666"
    "E36" = "  +12.16%  "
    "D37" = "'2.808"
    "E37" = "  +3.67%  "
    "E38" = "  +6.44%  "
    "D39" = "'2.957"
    "E39" = "  +4.07%  "
    "D40" = "'81.15"
    "E40" = "  +7.92%  "
    "D41" = "'6.698"
    "E41" = "  +7.62%  "
    "D42" = "'2.205"
    "E42" = "  +14.32%  "
    "D43" = "'0.4744"
    "E43" = "  +13.21%  "
    "D44" = "'0.8589"
    "E44" = "  +4.16%  "
    "D45" = "'104.91"
    "E45" = "  +4.62%  "
    "D46" = "'0.9986"
    "E46" = "  -0.02%  "
    "B47" = "EnergySwap"
    "C47" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D47" = "'9.982"
    "E47" = "  +4.21%  "
    "B48" = "Aptos"
    "C48" = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
    "D48" = "'7.667"
    "E48" = "  +10.15%  "
    "D49" = "'0.4373"
    "E49" = "  +12.94%  "
    "D50" = "'36.98"
    "E50" = "  +5.89%  "
    "B51" = "Maker"
    "C51" = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
    "D51" = "'933.14"
    "E51" = "  +5.91%  "
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $range = $ws.Range($addr)
    $range.Value = $value
    if ($value.StartsWith("'")) {
        # Drop the quote-prefix cell style Excel just applied so the cell
        # keeps its original (default) formatting - only the stored text
        # type needs to change, not the style index.
        $range.ClearFormats()
    }
}
